# Applies the "Fruta / hortaliza, semanal" update: the daily price rows for
# "Ramas de apio" (Agricola del Norte S.A. de Arica) get re-shuffled so that
# each row's Fecha/Variedad/Calidad/Volumen/Precio columns correspond to a
# different day's record than before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (row number => values for columns D,H,I,J,K,L,M,P)
# Column layout: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificacion.

$rows = @(
    @{ Row = 2; D = 44259; H = "Sin especificar"; I = "Primera";  J = 80;  K = 4000; L = 4500; M = 4250; P = 4250 },
    @{ Row = 3; D = 44263; H = "Sin especificar"; I = "Primera";  J = 100; K = 7000; L = 8000; M = 7500; P = 7500 },
    @{ Row = 4; D = 44253; H = "Americana (o)";   I = "Segunda";  J = 100; K = 4000; L = 4500; M = 4250; P = 4250 },
    @{ Row = 5; D = 44281; H = "Sin especificar"; I = "Primera";  J = 100; K = 5000; L = 6000; M = 5500; P = 5500 },
    @{ Row = 7; D = 44414; H = "Sin especificar"; I = "Primera";  J = 100; K = 6000; L = 7000; M = 6500; P = 6500 },
    @{ Row = 8; D = 44371; H = "Sin especificar"; I = "Primera";  J = 80;  K = 7000; L = 8000; M = 7375; P = 7375 },
    @{ Row = 9; D = 44410; H = "Sin especificar"; I = "Primera";  J = 100; K = 5500; L = 6000; M = 5750; P = 5750 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($row, 8).Value  = $r.H   # H - Variedad
    $ws.Cells.Item($row, 9).Value  = $r.I   # I - Calidad
    $ws.Cells.Item($row, 10).Value = $r.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio $/Kg
}
